# Adds "description for the results" sections (MapReduce problem / result
# narrative blocks) to Sheet1, Sheet2 and Sheet3 of the workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# Sheet3 ("Sheet3") -- bold the airline-name column, then add the two
# Mapreduce-problem / result call-out blocks below the data table.
# ---------------------------------------------------------------------
$ws3.Columns("A").Font.Bold = $true

$ws3.Range("A12").Value = "Mapreduce Problem:"
$ws3.Range("A12:B12").Merge()
$ws3.Range("A12:B12").Font.Bold = $true
$ws3.Range("A12:B12").Interior.Color = 65535
$ws3.Range("A12:B12").HorizontalAlignment = -4108

$ws3.Range("A13").Value = "Which airline has most negative feedback and which airline has most positive feedback"
$ws3.Range("A13:B13").Merge()
$ws3.Range("A13:B13").Font.Name = "Consolas"
$ws3.Range("A13:B13").Font.Size = 9
$ws3.Range("A13:B13").Font.Color = 2369838
$ws3.Range("A13:B13").Font.Family = 3
$ws3.Range("A13:B13").HorizontalAlignment = -4108
$ws3.Range("A13:B13").WrapText = $true
$ws3.Rows("13").RowHeight = 57.75

$ws3.Range("A15").Value = "Result"
$ws3.Range("A15:B15").Merge()
$ws3.Range("A15:B15").Font.Bold = $true
$ws3.Range("A15:B15").Interior.Color = 65535
$ws3.Range("A15:B15").HorizontalAlignment = -4108

$ws3.Range("A16").Value = "Based on the data united airlines has most negative feedback and Southwest airline has most positive feedback"
$ws3.Range("A16:B16").Merge()
$ws3.Range("A16:B16").HorizontalAlignment = -4108
$ws3.Range("A16:B16").WrapText = $true
$ws3.Rows("16").RowHeight = 73.5

$ws3.Range("A15").Select()

# ---------------------------------------------------------------------
# Sheet1 -- same kind of call-out blocks, different question/answer.
# ---------------------------------------------------------------------
$ws1.Range("A9").Value = "MapReduce Problem:"
$ws1.Range("A9:B9").Merge()
$ws1.Range("A9:B9").Font.Bold = $true
$ws1.Range("A9:B9").Interior.Color = 65535
$ws1.Range("A9:B9").HorizontalAlignment = -4108

$ws1.Range("A10").Value = "What is the positive feedback values for the United Airlines"
$ws1.Range("A10:B10").Merge()
$ws1.Range("A10:B10").Font.Name = "Consolas"
$ws1.Range("A10:B10").Font.Size = 9
$ws1.Range("A10:B10").Font.Color = 2369838
$ws1.Range("A10:B10").Font.Family = 3
$ws1.Range("A10:B10").HorizontalAlignment = -4108
$ws1.Range("A10:B10").WrapText = $true
$ws1.Rows("10").RowHeight = 39

$ws1.Range("A13").Value = "Result"
$ws1.Range("A13:B13").Merge()
$ws1.Range("A13:B13").Font.Bold = $true
$ws1.Range("A13:B13").Interior.Color = 65535
$ws1.Range("A13:B13").HorizontalAlignment = -4108

$ws1.Range("A14").Value = "Based on the above values we can say total 1218 is the count for negative and 238 is the positive feedback count"
$ws1.Range("A14:B14").Merge()
$ws1.Range("A14:B14").Font.Name = "Consolas"
$ws1.Range("A14:B14").Font.Size = 9
$ws1.Range("A14:B14").Font.Color = 2369838
$ws1.Range("A14:B14").Font.Family = 3
$ws1.Range("A14:B14").HorizontalAlignment = -4108
$ws1.Range("A14:B14").WrapText = $true
$ws1.Rows("14").RowHeight = 52.5

$ws1.Columns("A").ColumnWidth = 20.42578125

$ws1.Range("A9:B9").Select()

# ---------------------------------------------------------------------
# Sheet2 -- header renamed from "Reason" to "Reasons" (column A) /
# "Statistics for different negative feedbacks" (column B), plus two
# more Mapreduce-problem / result blocks appended below the table.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Reasons"
$ws2.Range("B1").Value = "Statistics for different negative feedbacks"

$ws2.Range("A14").Value = "MapReduce Problem:"
$ws2.Range("A14").Font.Bold = $true
$ws2.Range("A14").Interior.Color = 65535

$ws2.Range("A15").Value = "which type of issues made customers to give most negative feedback and its values"
$ws2.Range("A15").WrapText = $true
$ws2.Rows("15").RowHeight = 60

$ws2.Range("A18").Value = "Result:"
$ws2.Range("A18").Font.Bold = $true
$ws2.Range("A18").Interior.Color = 65535

$ws2.Range("A19").Value = "Customer service is the main issue for the neagtive feedback based on the resultant data"
$ws2.Range("A19").WrapText = $true
$ws2.Rows("19").RowHeight = 60

$ws2.Application.GoTo($ws2.Range("A14"))
